{"js": "// Generated replacements: each old expression is unique and maps to a unique new expression,\n// so a straightforward find-and-replace for each pair reproduces the diff exactly.\nconst replacements = [\n  [\"85+7=\", \"41+18=\"],\n  [\"50+4=\", \"55+33=\"],\n  [\"18+34=\", \"12+10=\"],\n  [\"5+52=\", \"56+37=\"],\n  [\"1+93=\", \"12+60=\"],\n  [\"22+75=\", \"22+45=\"],\n  [\"72+5=\", \"79-65=\"],\n  [\"94-21=\", \"66-49=\"],\n  [\"37+25=\", \"67+23=\"],\n  [\"8+4=\", \"22+69=\"],\n  [\"27+17=\", \"89-52=\"],\n  [\"16+28=\", \"75-68=\"],\n  [\"38-17=\", \"15+4=\"],\n  [\"25+9=\", \"7+18=\"],\n  [\"55+36=\", \"88-74=\"],\n  [\"73-7=\", \"11+60=\"],\n  [\"19+63=\", \"57-36=\"],\n  [\"99-21=\", \"33-30=\"],\n  [\"97-79=\", \"15+22=\"],\n  [\"26+60=\", \"2+94=\"],\n  [\"93+5=\", \"30+27=\"],\n  [\"70-22=\", \"56+40=\"],\n  [\"5+8=\", \"63+1=\"],\n  [\"92-28=\", \"82-57=\"],\n  [\"49-25=\", \"88-47=\"],\n  [\"22-5=\", \"66-19=\"],\n  [\"90-17=\", \"37+50=\"],\n  [\"82-76=\", \"39-3=\"],\n  [\"48-36=\", \"2+75=\"],\n  [\"97-71=\", \"11+4=\"],\n  [\"34+39=\", \"49-39=\"],\n  [\"14-7=\", \"88+1=\"],\n  [\"12+29=\", \"21+17=\"],\n  [\"64+32=\", \"61-47=\"],\n  [\"39-15=\", \"23-15=\"],\n  [\"34-20=\", \"70+28=\"],\n  [\"17+36=\", \"27-19=\"],\n  [\"82-38=\", \"74-36=\"],\n  [\"41-27=\", \"14+61=\"],\n  [\"29-14=\", \"95-47=\"],\n  [\"64-18=\", \"84-42=\"],\n  [\"40+7=\", \"73-14=\"],\n  [\"98-78=\", \"70-1=\"],\n  [\"0+52=\", \"97-86=\"],\n  [\"81-0=\", \"1+65=\"],\n  [\"67-11=\", \"96-30=\"],\n  [\"15+24=\", \"45-44=\"],\n  [\"66-13=\", \"50+38=\"],\n  [\"80+3=\", \"82-50=\"],\n  [\"7+60=\", \"68+15=\"],\n  [\"12-1=\", \"96-88=\"],\n  [\"94-31=\", \"79-20=\"],\n  [\"67-56=\", \"68-53=\"],\n  [\"72+18=\", \"62-7=\"],\n  [\"47-29=\", \"92-65=\"],\n  [\"88-29=\", \"76-65=\"],\n  [\"57+32=\", \"6-5=\"],\n  [\"32+41=\", \"85-1=\"],\n  [\"42+26=\", \"54-13=\"],\n  [\"62-35=\", \"18-17=\"],\n  [\"29-22=\", \"9+20=\"],\n  [\"55-12=\", \"36-24=\"],\n  [\"58-28=\", \"77-41=\"],\n  [\"20-0=\", \"12+21=\"],\n  [\"24+2=\", \"64-60=\"],\n  [\"16+82=\", \"70-13=\"],\n  [\"27-24=\", \"57-21=\"],\n  [\"77-66=\", \"6-0=\"],\n  [\"6+75=\", \"60-41=\"],\n  [\"12+15=\", \"93-44=\"],\n  [\"32+30=\", \"74-46=\"],\n  [\"67-8=\", \"67+26=\"],\n  [\"48-0=\", \"45+27=\"],\n  [\"83-76=\", \"61+38=\"],\n  [\"92-11=\", \"60-47=\"],\n  [\"3+19=\", \"54+5=\"],\n  [\"14+46=\", \"95-4=\"],\n  [\"96-18=\", \"8+23=\"],\n  [\"26+11=\", \"63-25=\"],\n  [\"35+31=\", \"53+1=\"],\n  [\"11+87=\", \"60-45=\"],\n  [\"86-26=\", \"0+40=\"],\n  [\"32-7=\", \"11+82=\"],\n  [\"16-0=\", \"59-24=\"],\n  [\"45-41=\", \"0+83=\"],\n  [\"65+4=\", \"97-36=\"],\n  [\"4+71=\", \"64-16=\"],\n  [\"53+17=\", \"17+70=\"],\n  [\"71+10=\", \"0+89=\"],\n  [\"85-32=\", \"84-7=\"],\n  [\"32+36=\", \"45-11=\"],\n  [\"99-79=\", \"98-33=\"],\n  [\"50+21=\", \"16-11=\"],\n  [\"35+46=\", \"90-60=\"],\n  [\"11+54=\", \"82+0=\"],\n  [\"36+1=\", \"18+70=\"],\n  [\"38+28=\", \"71+9=\"],\n  [\"16+1=\", \"84-69=\"],\n  [\"96-94=\", \"24+63=\"],\n  [\"97-40=\", \"0+51=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Generated replacements: each old expression is unique and maps to a unique new expression,\n# so a straightforward Find/Replace for each pair reproduces the diff exactly.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"85+7=\", \"41+18=\"),\n  @(\"50+4=\", \"55+33=\"),\n  @(\"18+34=\", \"12+10=\"),\n  @(\"5+52=\", \"56+37=\"),\n  @(\"1+93=\", \"12+60=\"),\n  @(\"22+75=\", \"22+45=\"),\n  @(\"72+5=\", \"79-65=\"),\n  @(\"94-21=\", \"66-49=\"),\n  @(\"37+25=\", \"67+23=\"),\n  @(\"8+4=\", \"22+69=\"),\n  @(\"27+17=\", \"89-52=\"),\n  @(\"16+28=\", \"75-68=\"),\n  @(\"38-17=\", \"15+4=\"),\n  @(\"25+9=\", \"7+18=\"),\n  @(\"55+36=\", \"88-74=\"),\n  @(\"73-7=\", \"11+60=\"),\n  @(\"19+63=\", \"57-36=\"),\n  @(\"99-21=\", \"33-30=\"),\n  @(\"97-79=\", \"15+22=\"),\n  @(\"26+60=\", \"2+94=\"),\n  @(\"93+5=\", \"30+27=\"),\n  @(\"70-22=\", \"56+40=\"),\n  @(\"5+8=\", \"63+1=\"),\n  @(\"92-28=\", \"82-57=\"),\n  @(\"49-25=\", \"88-47=\"),\n  @(\"22-5=\", \"66-19=\"),\n  @(\"90-17=\", \"37+50=\"),\n  @(\"82-76=\", \"39-3=\"),\n  @(\"48-36=\", \"2+75=\"),\n  @(\"97-71=\", \"11+4=\"),\n  @(\"34+39=\", \"49-39=\"),\n  @(\"14-7=\", \"88+1=\"),\n  @(\"12+29=\", \"21+17=\"),\n  @(\"64+32=\", \"61-47=\"),\n  @(\"39-15=\", \"23-15=\"),\n  @(\"34-20=\", \"70+28=\"),\n  @(\"17+36=\", \"27-19=\"),\n  @(\"82-38=\", \"74-36=\"),\n  @(\"41-27=\", \"14+61=\"),\n  @(\"29-14=\", \"95-47=\"),\n  @(\"64-18=\", \"84-42=\"),\n  @(\"40+7=\", \"73-14=\"),\n  @(\"98-78=\", \"70-1=\"),\n  @(\"0+52=\", \"97-86=\"),\n  @(\"81-0=\", \"1+65=\"),\n  @(\"67-11=\", \"96-30=\"),\n  @(\"15+24=\", \"45-44=\"),\n  @(\"66-13=\", \"50+38=\"),\n  @(\"80+3=\", \"82-50=\"),\n  @(\"7+60=\", \"68+15=\"),\n  @(\"12-1=\", \"96-88=\"),\n  @(\"94-31=\", \"79-20=\"),\n  @(\"67-56=\", \"68-53=\"),\n  @(\"72+18=\", \"62-7=\"),\n  @(\"47-29=\", \"92-65=\"),\n  @(\"88-29=\", \"76-65=\"),\n  @(\"57+32=\", \"6-5=\"),\n  @(\"32+41=\", \"85-1=\"),\n  @(\"42+26=\", \"54-13=\"),\n  @(\"62-35=\", \"18-17=\"),\n  @(\"29-22=\", \"9+20=\"),\n  @(\"55-12=\", \"36-24=\"),\n  @(\"58-28=\", \"77-41=\"),\n  @(\"20-0=\", \"12+21=\"),\n  @(\"24+2=\", \"64-60=\"),\n  @(\"16+82=\", \"70-13=\"),\n  @(\"27-24=\", \"57-21=\"),\n  @(\"77-66=\", \"6-0=\"),\n  @(\"6+75=\", \"60-41=\"),\n  @(\"12+15=\", \"93-44=\"),\n  @(\"32+30=\", \"74-46=\"),\n  @(\"67-8=\", \"67+26=\"),\n  @(\"48-0=\", \"45+27=\"),\n  @(\"83-76=\", \"61+38=\"),\n  @(\"92-11=\", \"60-47=\"),\n  @(\"3+19=\", \"54+5=\"),\n  @(\"14+46=\", \"95-4=\"),\n  @(\"96-18=\", \"8+23=\"),\n  @(\"26+11=\", \"63-25=\"),\n  @(\"35+31=\", \"53+1=\"),\n  @(\"11+87=\", \"60-45=\"),\n  @(\"86-26=\", \"0+40=\"),\n  @(\"32-7=\", \"11+82=\"),\n  @(\"16-0=\", \"59-24=\"),\n  @(\"45-41=\", \"0+83=\"),\n  @(\"65+4=\", \"97-36=\"),\n  @(\"4+71=\", \"64-16=\"),\n  @(\"53+17=\", \"17+70=\"),\n  @(\"71+10=\", \"0+89=\"),\n  @(\"85-32=\", \"84-7=\"),\n  @(\"32+36=\", \"45-11=\"),\n  @(\"99-79=\", \"98-33=\"),\n  @(\"50+21=\", \"16-11=\"),\n  @(\"35+46=\", \"90-60=\"),\n  @(\"11+54=\", \"82+0=\"),\n  @(\"36+1=\", \"18+70=\"),\n  @(\"38+28=\", \"71+9=\"),\n  @(\"16+1=\", \"84-69=\"),\n  @(\"96-94=\", \"24+63=\"),\n  @(\"97-40=\", \"0+51=\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}"}
